$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.376.08'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.567.59'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '290.83'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3769'
$ws.Range('E7').Value = '  +2.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.04'
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3396'
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07580'
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.136'
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.00'
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.958'
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.912'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.563.01'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001128'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.74'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06749'
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.58'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.202'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.93'
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.365.27'
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.705'
$ws.Range('E26').Value = '  -5.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.22'
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '147.95'
$ws.Range('E28').Value = '  +1.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.019'
$ws.Range('E29').Value = '  +0.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.64'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.740.28'
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9931'
$ws.Range('E33').Value = '  -2.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.056'
$ws.Range('E34').Value = '  -2.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '10.10'
$ws.Range('E35').Value = '  +0.55%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.436'
$ws.Range('E36').Value = '  +12.66%  '
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08458'
$ws.Range('E37').Value = '  -0.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02489'
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2291'
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06452'
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.423'
$ws.Range('E41').Value = '  -1.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6314'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.27'
$ws.Range('E43').Value = '  -3.45%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('E45').Value = '  -2.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.803'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5924'
$ws.Range('E47').Value = '  -0.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.073'
$ws.Range('E48').Value = '  -1.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.265'
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '124.67'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07324'
$ws.Range('E51').Value = '  +0.56%  '
